$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

# Re-trigger the merge border redistribution for the B1:D1 header band so that
# C1 and D1 pick up their own (top+bottom) / (right+top+bottom) thin border
# while B1 keeps its original full-box border.
$ws1.Range("B1:D1").UnMerge()
$ws1.Range("B1:D1").Merge()
$ws1.Range("B1").Borders.Item(10).LineStyle = 1

# Rename column header "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

# Same border redistribution trick for both merged header bands
$ws2.Range("B1:D1").UnMerge()
$ws2.Range("B1:D1").Merge()
$ws2.Range("B1").Borders.Item(10).LineStyle = 1

$ws2.Range("E1:G1").UnMerge()
$ws2.Range("E1:G1").Merge()
$ws2.Range("E1").Borders.Item(10).LineStyle = 1

# Rename column headers "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5
$ws2.Range("G5").ClearContents()
